# Append new closed/opened trade record (Trade #22) for base_strategy
# to both the "All Trades" log sheet and the per-strategy "base_strategy" sheet.

$wb = $excel.ActiveWorkbook

$tradeNumber   = 22
$tradeDate     = "2026-02-16"
$tradeTime     = "22:54:16"
$strategy      = "base_strategy"
$side          = "DOWN"
$entryPrice    = 49.999998
$exitPrice     = ""
$status        = "OPEN"
$pnlPct        = 0
$pnlDollar     = 0
$capitalAfter  = 100
$entrySlippage = 0
$exitSlippage  = 0
$confidence    = 0.6
$entryReason   = "Normal spread capture: 19600 bps"
$exitReason    = ""
$durationMin   = 0

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Next empty row is right after the current used range (row 23 here).
    $newRow = $ws.UsedRange.Rows.Count + 1

    $ws.Cells.Item($newRow, 1).Value = $tradeNumber

    # Prefix with an apostrophe so Excel stores the date as literal text
    # instead of auto-converting the "YYYY-MM-DD" pattern into a date serial.
    $ws.Cells.Item($newRow, 2).Value = "'" + $tradeDate

    $ws.Cells.Item($newRow, 3).Value  = $tradeTime
    $ws.Cells.Item($newRow, 4).Value  = $strategy
    $ws.Cells.Item($newRow, 5).Value  = $side
    $ws.Cells.Item($newRow, 6).Value  = $entryPrice
    $ws.Cells.Item($newRow, 7).Value  = $exitPrice
    $ws.Cells.Item($newRow, 8).Value  = $status
    $ws.Cells.Item($newRow, 9).Value  = $pnlPct
    $ws.Cells.Item($newRow, 10).Value = $pnlDollar
    $ws.Cells.Item($newRow, 11).Value = $capitalAfter
    $ws.Cells.Item($newRow, 12).Value = $entrySlippage
    $ws.Cells.Item($newRow, 13).Value = $exitSlippage
    $ws.Cells.Item($newRow, 14).Value = $confidence
    $ws.Cells.Item($newRow, 15).Value = $entryReason
    $ws.Cells.Item($newRow, 16).Value = $exitReason
    $ws.Cells.Item($newRow, 17).Value = $durationMin

    Write-Host "Updated '$sheetName': wrote trade #$tradeNumber to row $newRow"
}
